# Auto-generated edit script applying the Asura_Profits.xlsx market-data refresh
# described in the commit diff. All changes are plain numeric value updates
# (cached market-board figures), including a few cells that are newly added
# or fully cleared in the target state.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 308.13043
$ws.Range("I33").Value = 308.8095
$ws.Range("K33").Value = 308.8095
$ws.Range("M33").Value = -79.80950000000001
$ws.Range("H52").Value = 9000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H53").Value = 32.3
$ws.Range("J53").Value = 49.5
$ws.Range("L53").Value = 49.5
$ws.Range("N53").Value = -1323.5
$ws.Range("H62").Value = 4101.5
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 3802
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 3802
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -5050
$ws.Range("H65").Value = 4101.5
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 3802
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 19010
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -25250
$ws.Range("H69").Value = 4248.75
$ws.Range("I69").Value = 2495
$ws.Range("J69").Value = 4833.3335
$ws.Range("K69").Value = 7485
$ws.Range("L69").Value = 14500.0005
$ws.Range("M69").Value = -6611
$ws.Range("N69").Value = -16248.0005
$ws.Range("H72").Value = 4248.75
$ws.Range("I72").Value = 2495
$ws.Range("J72").Value = 4833.3335
$ws.Range("K72").Value = 22455
$ws.Range("L72").Value = 43500.0015
$ws.Range("M72").Value = -18087
$ws.Range("N72").Value = -52236.0015
$ws.Range("H116").Value = 25003688
$ws.Range("I116").Value = 33337334
$ws.Range("J116").Value = 2747.5
$ws.Range("K116").Value = 33337334
$ws.Range("L116").Value = 2747.5
$ws.Range("M116").Value = -33333892
$ws.Range("N116").Value = -9631.5
$ws.Range("H140").Value = 104554.445
$ws.Range("J140").Value = 105930
$ws.Range("L140").Value = 105930
$ws.Range("N140").Value = -116290

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16025.019
$ws.Range("I32").Value = 14338.24
$ws.Range("J32").Value = 44138
$ws.Range("K32").Value = 14338.24
$ws.Range("L32").Value = 44138
$ws.Range("M32").Value = -14051.24
$ws.Range("N32").Value = -44712
$ws.Range("H45").Value = 957.38464
$ws.Range("I45").Value = 804
$ws.Range("K45").Value = 804
$ws.Range("M45").Value = -427
$ws.Range("H60").Value = 50000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H74").Value = 1203.0857
$ws.Range("I74").Value = 1081.5
$ws.Range("K74").Value = 1081.5
$ws.Range("M74").Value = -207.5
$ws.Range("H77").Value = 1203.0857
$ws.Range("I77").Value = 1081.5
$ws.Range("K77").Value = 5407.5
$ws.Range("M77").Value = -1039.5
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -504
$ws.Range("H102").Value = 2900
$ws.Range("I102").Value = 1516.6666
$ws.Range("K102").Value = 1516.6666
$ws.Range("M102").Value = 105.3334
$ws.Range("H132").Value = 477862.2
$ws.Range("I132").Value = 607261.2
$ws.Range("J132").Value = 3399.111
$ws.Range("K132").Value = 1821783.6
$ws.Range("L132").Value = 10197.333
$ws.Range("M132").Value = -1819253.6
$ws.Range("N132").Value = -15257.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 231666.5
$ws.Range("J43").Value = 231666.5
$ws.Range("L43").Value = 231666.5
$ws.Range("N43").Value = -232028.5
$ws.Range("H99").Value = 1182.5
$ws.Range("I99").Value = 918
$ws.Range("K99").Value = 918
$ws.Range("M99").Value = 580
$ws.Range("H105").Value = 3667.5625
$ws.Range("I105").Value = 3245.4
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 3245.4
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -1498.4
$ws.Range("N105").Value = -13494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 30000
$ws.Range("K33").Value = 30000
$ws.Range("M33").Value = -29621
$ws.Range("H132").Value = 1796.4572
$ws.Range("I132").Value = 1499.3125
$ws.Range("K132").Value = 4497.9375
$ws.Range("M132").Value = -1967.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7194.875
$ws.Range("J75").Value = 7194.875
$ws.Range("L75").Value = 21584.625
$ws.Range("N75").Value = -23580.625
$ws.Range("H78").Value = 7194.875
$ws.Range("J78").Value = 7194.875
$ws.Range("L78").Value = 64753.875
$ws.Range("N78").Value = -74737.875
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 26750
$ws.Range("J102").Value = 3500
$ws.Range("L102").Value = 10500
$ws.Range("N102").Value = -15368
$ws.Range("H110").Value = 9900
$ws.Range("I110").Value = 6500
$ws.Range("K110").Value = 19500
$ws.Range("M110").Value = -15410
$ws.Range("H119").Value = 4994
$ws.Range("I119").Value = 2991.6
$ws.Range("K119").Value = 8974.799999999999
$ws.Range("M119").Value = -4136.799999999999
$ws.Range("H120").Value = 11321.667
$ws.Range("J120").Value = 13633.333
$ws.Range("L120").Value = 40899.999
$ws.Range("N120").Value = -50575.999
$ws.Range("H122").Value = 3274.7
$ws.Range("J122").Value = 6074.9
$ws.Range("L122").Value = 54674.1
$ws.Range("N122").Value = -59574.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 36138
$ws.Range("I97").Value = 40928.46
$ws.Range("K97").Value = 40928.46
$ws.Range("M97").Value = -40432.46
$ws.Range("H132").Value = 2206.6086
$ws.Range("I132").Value = 1368.8462
$ws.Range("J132").Value = 3295.7
$ws.Range("K132").Value = 4106.5386
$ws.Range("L132").Value = 9887.099999999999
$ws.Range("M132").Value = -1576.5386
$ws.Range("N132").Value = -14947.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1928.5714
$ws.Range("I68").Value = 1700
$ws.Range("K68").Value = 1700
$ws.Range("M68").Value = -951
$ws.Range("H71").Value = 1928.5714
$ws.Range("I71").Value = 1700
$ws.Range("K71").Value = 8500
$ws.Range("M71").Value = -4756
$ws.Range("H100").Value = 12580
$ws.Range("I100").Value = 18966.666
$ws.Range("K100").Value = 18966.666
$ws.Range("M100").Value = -18425.666
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15626287
$ws.Range("I122").Value = 27778866
$ws.Range("J122").Value = 1542.5714
$ws.Range("K122").Value = 83336598
$ws.Range("L122").Value = 4627.7142
$ws.Range("M122").Value = -83334148
$ws.Range("N122").Value = -9527.7142
